$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Oct 25 12:07:08 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 12:07:21 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 12:07:34 EDT 2024"
